# Update the "F" column (sales/ticket count) numbers on the "展览" and
# "全部类型" sheets (and the single matching row on "演出") to reflect the
# newly scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 778
    4  = 1565
    6  = 104
    7  = 180
    8  = 2
    9  = 6372
    11 = 414
    13 = 5549
    16 = 1215
    18 = 69
    20 = 77
    22 = 316
    23 = 32
    25 = 3980
    26 = 12
    27 = 176
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (sheet2) ----------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 102

# --- Sheet "全部类型" (sheet4) -------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 102
    4  = 778
    5  = 1565
    7  = 104
    8  = 180
    9  = 2
    10 = 6372
    12 = 414
    14 = 5550
    17 = 1215
    19 = 69
    21 = 77
    23 = 316
    24 = 32
    26 = 3980
    28 = 12
    29 = 176
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
